$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I1").Value2 = "on_topic"
$ws.Range("I1").Font.Name = "Arial"
$ws.Range("I1").Font.Size = 10
Write-Host "done"
